$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need an explicit Text
# number format first, otherwise Excel auto-converts the assigned string
# into a numeric value (losing the original text representation, e.g.
# trailing zeros such as "15.10" -> 15.1).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.603.65"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.508.21"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "315.73"
$ws.Range("E5").Value = "  +4.14%  "
$ws.Range("D6").Value = "94.91"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  -1.95%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "35.73"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "7.55"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Value = "2.891.92"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.504.84"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.10"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "0.847"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "42.680.47"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  +6.11%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "0.0₃0958"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "69.46"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "250.75"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "2.94"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").Value = "2.07"
$ws.Range("E25").Value = "  +3.02%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "26.29"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "2.42"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("D29").Value = "41.39"
$ws.Range("E29").Value = "  +11.71%  "
$ws.Range("D30").Value = "10.25"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "5.94"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").Value = "158.99"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").Value = "2.12"
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.67"
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "3.25"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").Value = "0.0777"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "23.77"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  +17.58%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0305"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "3.76"
$ws.Range("E45").Value = "  -1.58%  "
$ws.Range("D46").Value = "2.022.45"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "84.86"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "8.90"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "74.34"
$ws.Range("E49").Value = "  +3.11%  "
$ws.Range("D50").Value = "2.749.65"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "102.04"
$ws.Range("E51").Value = "  +1.55%  "
